# Production shares update MG 0206
# Insert a new "Water electrolysis" row above the existing "CO$_2$ utilization"
# row (currently row 8), pushing that row and the ones below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8; rows 8-11 (CO$_2$ utilization ... Plastic
# waste recycling with CC) shift down to rows 9-12.
$ws.Rows(8).Insert()

# Fill in the new "Water electrolysis" row with zero values across all years.
$ws.Range("A8").Value = "Water electrolysis"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Match the formatting (bold, centered, thin-bordered label style) used by
# the other row-label cells in column A.
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$ws.Range("A8").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$ws.Range("A8").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("A8").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
